$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date serial from 45184 to 45186 for every data row (2-199)
for ($row = 2; $row -le 199; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# Add the friendly display-text argument (the case/beteckning id from column A)
# to the HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2 and 3.
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

foreach ($row in 2, 3) {
    $label = $ws.Range("A$row").Value2
    foreach ($col in $hyperlinkCols) {
        $cellRef = "$col$row"
        $formula = $ws.Range($cellRef).Formula
        if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $Matches[1]
            $ws.Range($cellRef).Formula = "=HYPERLINK(`"$url`", `"$label`")"
        }
    }
}
